# Auto-generated Excel COM-interop script that applies the numeric updates
# described by the "chore: update Sheets via scheduled runner" diff.
# Each Leve row below has one or more of its price/profit columns (H-N)
# refreshed to new Market Board snapshot values.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 5854668
$ws.Range("I51").Value = 13892124
$ws.Range("J51").Value = 9245.454
$ws.Range("K51").Value = 13892124
$ws.Range("L51").Value = 9245.454
$ws.Range("M51").Value = -13891640
$ws.Range("N51").Value = -10213.454
# Row 96
$ws.Range("H96").Value = 58834640
$ws.Range("I96").Value = 4040.5
$ws.Range("K96").Value = 12121.5
$ws.Range("M96").Value = -10748.5
# Row 113
$ws.Range("H113").Value = 2184.6155
$ws.Range("I113").Value = 1950
$ws.Range("J113").Value = 2227.2727
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 2227.2727
$ws.Range("M113").Value = 1304
$ws.Range("N113").Value = -8735.2727
# Row 129
$ws.Range("H129").Value = 1058.3636
$ws.Range("I129").Value = 5298.5
$ws.Range("J129").Value = 970.9382000000001
$ws.Range("K129").Value = 15895.5
$ws.Range("L129").Value = 2912.8146
$ws.Range("M129").Value = -10895.5
$ws.Range("N129").Value = -12912.8146
# Row 137
$ws.Range("H137").Value = 1712756.5
$ws.Range("I137").Value = 2849897.8
$ws.Range("J137").Value = 7044.722
$ws.Range("K137").Value = 8549693.399999999
$ws.Range("L137").Value = 21134.166
$ws.Range("M137").Value = -8547143.399999999
$ws.Range("N137").Value = -26234.166

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11796.667
$ws.Range("I32").Value = 10979.276
$ws.Range("J32").Value = 17284.857
$ws.Range("K32").Value = 10979.276
$ws.Range("L32").Value = 17284.857
$ws.Range("M32").Value = -10692.276
$ws.Range("N32").Value = -17858.857
# Row 61
$ws.Range("H61").Value = 2767.375
$ws.Range("I61").Value = 1990.4546
$ws.Range("K61").Value = 1990.4546
$ws.Range("M61").Value = -1778.4546
# Row 102
$ws.Range("H102").Value = 52646.223
$ws.Range("J102").Value = 58752.125
$ws.Range("L102").Value = 58752.125
$ws.Range("N102").Value = -61996.125
# Row 110
$ws.Range("H110").Value = 1437.4482
$ws.Range("I110").Value = 1372.7142
$ws.Range("J110").Value = 1607.375
$ws.Range("K110").Value = 1372.7142
$ws.Range("L110").Value = 1607.375
$ws.Range("M110").Value = 672.2858000000001
$ws.Range("N110").Value = -5697.375
# Row 136
$ws.Range("H136").Value = 2767.375
$ws.Range("I136").Value = 1990.4546
$ws.Range("K136").Value = 5971.3638
$ws.Range("M136").Value = -3421.3638

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2085.8667
$ws.Range("I20").Value = 1248.6
$ws.Range("J20").Value = 2504.5
$ws.Range("K20").Value = 1248.6
$ws.Range("L20").Value = 2504.5
$ws.Range("M20").Value = -1001.6
$ws.Range("N20").Value = -2998.5
# Row 94
$ws.Range("H94").Value = 620.3913
$ws.Range("I94").Value = 623.36584
$ws.Range("J94").Value = 596
$ws.Range("K94").Value = 623.36584
$ws.Range("L94").Value = 596
$ws.Range("M94").Value = -172.36584
$ws.Range("N94").Value = -1498
# Row 105
$ws.Range("H105").Value = 2335.5881
$ws.Range("I105").Value = 1792.2106
$ws.Range("J105").Value = 3023.8667
$ws.Range("K105").Value = 1792.2106
$ws.Range("L105").Value = 3023.8667
$ws.Range("M105").Value = -45.21060000000011
$ws.Range("N105").Value = -6517.8667
# Row 107
$ws.Range("H107").Value = 1760.125
$ws.Range("I107").Value = 1355.35
$ws.Range("K107").Value = 1355.35
$ws.Range("M107").Value = 564.6500000000001

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1040
$ws.Range("I16").Value = 1076
$ws.Range("K16").Value = 1076
$ws.Range("M16").Value = -789
# Row 22
$ws.Range("H22").Value = 1657.8125
$ws.Range("J22").Value = 5525
$ws.Range("L22").Value = 5525
$ws.Range("N22").Value = -6225
# Row 99
$ws.Range("H99").Value = 2787.2856
$ws.Range("I99").Value = 2270.6667
$ws.Range("J99").Value = 3174.75
$ws.Range("K99").Value = 2270.6667
$ws.Range("L99").Value = 3174.75
$ws.Range("M99").Value = -772.6667000000002
$ws.Range("N99").Value = -6170.75
# Row 113
$ws.Range("H113").Value = 1040
$ws.Range("I113").Value = 1076
$ws.Range("K113").Value = 1076
$ws.Range("M113").Value = 1094
# Row 126
$ws.Range("H126").Value = 2787.2856
$ws.Range("I126").Value = 2270.6667
$ws.Range("J126").Value = 3174.75
$ws.Range("K126").Value = 6812.000100000001
$ws.Range("L126").Value = 9524.25
$ws.Range("M126").Value = -4342.000100000001
$ws.Range("N126").Value = -14464.25
# Row 132
$ws.Range("H132").Value = 53475.703
$ws.Range("I132").Value = 1107.45
$ws.Range("J132").Value = 203099.28
$ws.Range("K132").Value = 3322.35
$ws.Range("L132").Value = 609297.84
$ws.Range("M132").Value = -792.3500000000004
$ws.Range("N132").Value = -614357.84
# Row 134
$ws.Range("H134").Value = 781577
$ws.Range("I134").Value = 504058.78
$ws.Range("J134").Value = 1752890.6
$ws.Range("K134").Value = 1512176.34
$ws.Range("L134").Value = 5258671.800000001
$ws.Range("M134").Value = -1509641.34
$ws.Range("N134").Value = -5263741.800000001

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 3504.389
$ws.Range("J5").Value = 1429
$ws.Range("L5").Value = 4287
$ws.Range("N5").Value = -4511
# Row 112
$ws.Range("H112").Value = 16669079
$ws.Range("J112").Value = 3035.5557
$ws.Range("L112").Value = 9106.667099999999
$ws.Range("N112").Value = -11322.6671
# Row 135
$ws.Range("H135").Value = 3504.389
$ws.Range("J135").Value = 1429
$ws.Range("L135").Value = 12861
$ws.Range("N135").Value = -17931

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 2923.2354
$ws.Range("I97").Value = 1892.2222
$ws.Range("J97").Value = 6900
$ws.Range("K97").Value = 1892.2222
$ws.Range("L97").Value = 6900
$ws.Range("M97").Value = -1396.2222
$ws.Range("N97").Value = -7892
# Row 113
$ws.Range("H113").Value = 1470.7059
$ws.Range("J113").Value = 1731.3334
$ws.Range("L113").Value = 1731.3334
$ws.Range("N113").Value = -6071.3334
# Row 122
$ws.Range("H122").Value = 1434.7693
$ws.Range("I122").Value = 1064.5714
$ws.Range("J122").Value = 1866.6666
$ws.Range("K122").Value = 3193.7142
$ws.Range("L122").Value = 5599.9998
$ws.Range("M122").Value = -743.7142000000003
$ws.Range("N122").Value = -10499.9998
# Row 126
$ws.Range("H126").Value = 16815.857
$ws.Range("I126").Value = 22902.2
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 68706.60000000001
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = -66236.60000000001
$ws.Range("N126").Value = -9740
# Row 132
$ws.Range("H132").Value = 22224648
$ws.Range("I132").Value = 31251472
$ws.Range("J132").Value = 4769.6924
$ws.Range("K132").Value = 93754416
$ws.Range("L132").Value = 14309.0772
$ws.Range("M132").Value = -93751886
$ws.Range("N132").Value = -19369.0772
# Row 135
$ws.Range("H135").Value = 58499.5
$ws.Range("J135").Value = 58499.5
$ws.Range("L135").Value = 58499.5
$ws.Range("N135").Value = -68639.5
# Row 138
$ws.Range("H138").Value = 53500
$ws.Range("J138").Value = 53500
$ws.Range("L138").Value = 53500
$ws.Range("N138").Value = -63780
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
# Row 141
$ws.Range("H141").Value = 76142.664
$ws.Range("J141").Value = 76142.664
$ws.Range("L141").Value = 76142.664
$ws.Range("N141").Value = -86502.664

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3466.6667
$ws.Range("I40").Value = 2168.75
$ws.Range("K40").Value = 2168.75
$ws.Range("M40").Value = -2032.75
# Row 93
$ws.Range("H93").Value = 1345.2307
$ws.Range("I93").Value = 914.3333
$ws.Range("J93").Value = 1714.5714
$ws.Range("K93").Value = 914.3333
$ws.Range("L93").Value = 1714.5714
$ws.Range("M93").Value = 333.6667
$ws.Range("N93").Value = -4210.5714
# Row 132
$ws.Range("H132").Value = 2478.9834
$ws.Range("I132").Value = 1959.3024
$ws.Range("J132").Value = 3793.4707
$ws.Range("K132").Value = 5877.9072
$ws.Range("L132").Value = 11380.4121
$ws.Range("M132").Value = -3347.9072
$ws.Range("N132").Value = -16440.4121

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1100
$ws.Range("I96").Value = 1100
$ws.Range("K96").Value = 1100
$ws.Range("M96").Value = 273
# Row 100
$ws.Range("H100").Value = 360.125
$ws.Range("I100").Value = 220
$ws.Range("J100").Value = 444.2
$ws.Range("K100").Value = 440
$ws.Range("L100").Value = 888.4
$ws.Range("M100").Value = 101
$ws.Range("N100").Value = -1970.4
# Row 132
$ws.Range("H132").Value = 1403912.4
$ws.Range("I132").Value = 1611395.6
$ws.Range("J132").Value = 3399.75
$ws.Range("K132").Value = 4834186.800000001
$ws.Range("L132").Value = 10199.25
$ws.Range("M132").Value = -4831656.800000001
$ws.Range("N132").Value = -15259.25
# Row 135
$ws.Range("H135").Value = 31351.428
$ws.Range("J135").Value = 31351.428
$ws.Range("L135").Value = 31351.428
$ws.Range("N135").Value = -41491.428
# Row 137
$ws.Range("H137").Value = 62107.5
$ws.Range("J137").Value = 62107.5
$ws.Range("L137").Value = 62107.5
$ws.Range("N137").Value = -72307.5
# Row 139
$ws.Range("H139").Value = 45204
$ws.Range("J139").Value = 45204
$ws.Range("L139").Value = 45204
$ws.Range("N139").Value = -55484
# Row 141
$ws.Range("H141").Value = 36182
$ws.Range("J141").Value = 36182
$ws.Range("L141").Value = 36182
$ws.Range("N141").Value = -46542

